# Concentrado_Mejoras_2015.xlsx — "Se edito el nombre de la pestaña. Actualización menor"
#
# The commit renames the single worksheet tab from "Catalogo de cursos" to
# "Concentrado de cursos" and leaves a small trace of the editing session in
# the saved view state (the cell that was selected when the file was saved
# moved from E11 to B11).
#
# Notes on scope: the source diff also shows a couple of purely-cosmetic,
# non-semantic artifacts of a real Excel desktop save (the absolute path
# breadcrumb under mc:AlternateContent, and a cellXfs re-shuffle that keeps
# every cell's visible formatting identical, just renumbered). Those aren't
# things a user action maps to — they're host/save-session bookkeeping that
# isn't exposed through the Excel object model, so they're intentionally not
# reproduced here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "Concentrado de cursos"

# Leave the workbook's saved selection on B11 (was E11).
$ws.Range("B11").Select()
